$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ActivityCompleteBookingFlows")

# Update the test scenario row (row 2) with new Activity scenario data:
# EndDate changes from 22 to 16
$ws.Range("D2").Value = 16

# ShortLocation / Location change from "Las Vegas, NV, US" to "Los Angeles, CA, US"
$ws.Range("E2").Value = "Los Angeles, CA, US"
$ws.Range("F2").Value = "Los Angeles, CA, US"

# Row grew taller to fit the wrapped new text
$ws.Rows.Item(2).RowHeight = 31.5

# Reset the sheet selection back to the default top-left cell
$ws.Range("A1").Select()
